$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-04 Tuesday", "2024-06-05 Wednesday"),
    @("702÷8=87, 6", "390÷3=130, 0"),
    @("266÷2=133, 0", "625÷3=208, 1"),
    @("588÷3=196, 0", "925÷4=231, 1"),
    @("385÷2=192, 1", "143÷2=71, 1"),
    @("987÷6=164, 3", "509÷3=169, 2"),
    @("641÷9=71, 2", "741÷3=247, 0"),
    @("383÷9=42, 5", "724÷7=103, 3"),
    @("196÷9=21, 7", "885÷7=126, 3"),
    @("498÷4=124, 2", "925÷7=132, 1"),
    @("950÷7=135, 5", "823÷5=164, 3"),
    @("922÷3=307, 1", "371÷9=41, 2"),
    @("415÷7=59, 2", "471÷9=52, 3"),
    @("931÷5=186, 1", "139÷2=69, 1"),
    @("163÷4=40, 3", "484÷3=161, 1"),
    @("180÷5=36, 0", "228÷9=25, 3"),
    @("135÷2=67, 1", "796÷5=159, 1"),
    @("205÷4=51, 1", "893÷3=297, 2"),
    @("233÷6=38, 5", "123÷4=30, 3"),
    @("205÷9=22, 7", "743÷6=123, 5"),
    @("195÷9=21, 6", "398÷2=199, 0"),
    @("592÷3=197, 1", "106÷2=53, 0"),
    @("921÷6=153, 3", "562÷9=62, 4"),
    @("907÷4=226, 3", "624÷2=312, 0"),
    @("797÷6=132, 5", "367÷5=73, 2"),
    @("791÷9=87, 8", "497÷8=62, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
